# Pairwise sensitivity models added.
#
# Populates the "insulin_resistant_y_n" column (Z) for every study row that
# already has PCOS-group insulin-resistance data reported in column Y, using
# "y"/"n" to flag whether that row's cohort was classified as insulin
# resistant. Rows with no corresponding data (8, 10-14, 17-25, 28, 29) are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ZCOL = 26   # column Z

# row -> "y"/"n" for column Z (insulin_resistant_y_n)
$values = @{
    2  = "y"
    3  = "n"
    4  = "y"
    5  = "n"
    6  = "y"
    7  = "n"
    9  = "y"
    15 = "y"
    16 = "y"
    26 = "y"
    27 = "n"
    30 = "n"
    31 = "y"
    32 = "y"
    33 = "y"
    34 = "n"
    35 = "y"
    36 = "y"
    37 = "y"
    38 = "y"
    39 = "y"
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, $ZCOL).Value = $values[$row]
}

# Match the author's final view/selection state (scrolled to show column T
# onward, frozen pane scrolled to row 29, last touched cell Z39 selected).
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 20
$win.ScrollRow = 29
$ws.Range("Z39").Select() | Out-Null
